# Week 3 assignment submitted.
# Fill in the "Actual time length to complete" (column C) values for the
# tasks that were just finished, on the "week3" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week3")
$ws.Activate()

$ws.Range("C11").Value = 0.010416666666666666
$ws.Range("C12").Value = 0.013888888888888888
$ws.Range("C16").Value = 0.010416666666666666
$ws.Range("C17").Value = 0.013888888888888888

# Move / leave the active selection on C13, matching the saved view state.
$ws.Range("C13").Select()
